$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column F header placeholder (blank cell, same look as the other headers) ---
$ws.Range("F1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("F1").VerticalAlignment = -4108     # xlCenter
$ws.Range("F1").Font.Bold = $true

# --- Format the new derived-error columns (G:I) before filling them in, so the ---
# --- number format / alignment settle into a single style before formulas land ---
$calcRange = $ws.Range("G2:I22")
$calcRange.HorizontalAlignment = -4108        # xlCenter
$calcRange.VerticalAlignment = -4108          # xlCenter
$calcRange.NumberFormat = "0.000"

for ($r = 2; $r -le 22; $r++) {
    $ws.Cells.Item($r, 7).Formula = "=A$r*(1+(B$r/100))"
    $ws.Cells.Item($r, 8).Formula = "=A$r*(1+(C$r/100))"
    $ws.Cells.Item($r, 9).Formula = "=A$r*(1+(D$r/100))"
}

$ws.Range("F1").Select()
